$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.529.63'
$ws.Range('E2').Value = '  +2.84%  '
$ws.Range('D3').Value = '2.432.72'
$ws.Range('E3').Value = '  +0.23%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.60'
$ws.Range('E5').Value = '  +2.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.37'
$ws.Range('E6').Value = '  +4.38%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.514'
$ws.Range('E8').Value = '  +1.30%  '
$ws.Range('E9').Value = '  +6.65%  '
$ws.Range('D10').Value = '2.431.68'
$ws.Range('E10').Value = '  +0.22%  '
$ws.Range('E11').Value = '  -2.05%  '
$ws.Range('E12').Value = '  +2.03%  '
$ws.Range('E13').Value = '  -2.25%  '
$ws.Range('E14').Value = '  +4.77%  '
$ws.Range('D15').Value = '69.442.01'
$ws.Range('E15').Value = '  +2.87%  '
$ws.Range('D16').Value = '2.881.52'
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '23.96'
$ws.Range('E17').Value = '  +4.44%  '
$ws.Range('D18').Value = '2.427.72'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.80'
$ws.Range('E19').Value = '  +4.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '341.99'
$ws.Range('E20').Value = '  +3.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.15'
$ws.Range('E21').Value = '  +4.48%  '
$ws.Range('E22').Value = '  +2.85%  '
$ws.Range('E23').Value = '  +6.95%  '
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.90'
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').Value = '  +5.85%  '
$ws.Range('E27').Value = '  +5.63%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.995'
$ws.Range('E29').Value = '  -2.02%  '
$ws.Range('D30').Value = '0.0₃0850'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.38'
$ws.Range('E31').Value = '  +5.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.23'
$ws.Range('E32').Value = '  +10.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '452.87'
$ws.Range('E33').Value = '  +9.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('E35').Value = '  +1.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '157.76'
$ws.Range('E36').Value = '  -1.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.12'
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('E38').Value = '  +5.69%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.24'
$ws.Range('E40').Value = '  +2.57%  '
$ws.Range('E41').Value = '  +2.80%  '
$ws.Range('E42').Value = '  +3.80%  '
$ws.Range('E43').Value = '  +4.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '37.96'
$ws.Range('E44').Value = '  +1.55%  '
$ws.Range('E45').Value = '  +2.05%  '
$ws.Range('E46').Value = '  +5.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '134.99'
$ws.Range('E47').Value = '  +3.94%  '
$ws.Range('E48').Value = '  +2.69%  '
$ws.Range('E49').Value = '  +2.67%  '
$ws.Range('E50').Value = '  +2.79%  '
$ws.Range('E51').Value = '  +1.64%  '
